$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain number-looking string must be forced
# to Text (matching the source data, which stores all Price/Volume figures as
# text) so Excel does not silently reinterpret them as numeric values.

# Row 2
$ws.Range("D2").Value = "63.146.59"
$ws.Range("E2").Value = "  -1.56%  "

# Row 3
$ws.Range("D3").Value = "3.228.80"
$ws.Range("E3").Value = "  -1.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.53"
$ws.Range("E5").Value = "  +1.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.98"
$ws.Range("E6").Value = "  -4.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "3.227.13"
$ws.Range("E9").Value = "  -1.82%  "

# Row 10
$ws.Range("E10").Value = "  -2.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.02"
$ws.Range("E11").Value = "  -7.60%  "

# Row 12
$ws.Range("E12").Value = "  +2.00%  "

# Row 13
$ws.Range("E13").Value = "  -0.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.10"
$ws.Range("E14").Value = "  +0.58%  "

# Row 15
$ws.Range("D15").Value = "3.745.78"
$ws.Range("E15").Value = "  -1.43%  "

# Row 16
$ws.Range("E16").Value = "  -2.20%  "

# Row 17
$ws.Range("D17").Value = "3.235.03"
$ws.Range("E17").Value = "  -1.24%  "

# Row 18
$ws.Range("D18").Value = "63.000.95"
$ws.Range("E18").Value = "  -1.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.14"
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  +2.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.964"
$ws.Range("E21").Value = "  +2.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "366.06"
$ws.Range("E22").Value = "  -1.41%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.75"
$ws.Range("E23").Value = "  +3.40%  "

# Row 24
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").Value = "  +3.64%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.11"
$ws.Range("E25").Value = "  +1.41%  "

# Row 26
$ws.Range("E26").Value = "  +6.01%  "

# Row 27
$ws.Range("E27").Value = "  +1.74%  "

# Row 28
$ws.Range("E28").Value = "  -0.01%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.26"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  -1.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.45"
$ws.Range("E31").Value = "  +0.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "638.79"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.44"
$ws.Range("E33").Value = "  -3.25%  "

# Row 34
$ws.Range("E34").Value = "  +1.47%  "

# Row 35
$ws.Range("E35").Value = "  +2.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.95"
$ws.Range("E36").Value = "  -3.64%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.27%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.60"
$ws.Range("E38").Value = "  +2.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.376"
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.46%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0712"
$ws.Range("E41").Value = "  +8.84%  "

# Row 42
$ws.Range("E42").Value = "  +12.24%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.123"
$ws.Range("E43").Value = "  +1.36%  "

# Row 44
$ws.Range("D44").Value = "2.879.91"
$ws.Range("E44").Value = "  +1.19%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.94"
$ws.Range("E45").Value = "  +7.44%  "

# Row 46
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.69"
$ws.Range("E46").Value = "  +3.73%  "

# Row 47
$ws.Range("E47").Value = "  +3.24%  "

# Row 48
$ws.Range("E48").Value = "  +4.52%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("E49").Value = "  -2.76%  "

# Row 50
$ws.Range("E50").Value = "  +1.00%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.28"
$ws.Range("E51").Value = "  +1.36%  "
